$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '48.041.00'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.494.05'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.31'
$ws.Range('E5').Value = '  -1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.52'
$ws.Range('E6').Value = '  -3.06%  '
$ws.Range('E7').Value = '  -1.40%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  -4.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.75'
$ws.Range('E10').Value = '  -3.96%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  -2.16%  '
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.08'
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.886.14'
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.498.90'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.829'
$ws.Range('E17').Value = '  -3.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.890.32'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('B19').Value = 'ImmutableX'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.97'
$ws.Range('E19').Value = '  +9.08%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.95'
$ws.Range('E20').Value = '  -2.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.62'
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  -1.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.98'
$ws.Range('E23').Value = '  -1.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.68'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.73'
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.28'
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.70'
$ws.Range('E29').Value = '  -4.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.139'
$ws.Range('E30').Value = '  -3.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.67'
$ws.Range('E31').Value = '  -0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.29'
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.96'
$ws.Range('E34').Value = '  -5.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.27'
$ws.Range('E35').Value = '  -2.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0770'
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.57'
$ws.Range('E38').Value = '  -3.10%  '
$ws.Range('E39').Value = '  -4.38%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '121.64'
$ws.Range('E40').Value = '  +2.51%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.23'
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.110'
$ws.Range('E42').Value = '  -2.28%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.16'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.994.91'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.86'
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.14'
$ws.Range('E50').Value = '  -1.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.80'
$ws.Range('E51').Value = '  -1.93%  '
